# Fix probe mapping to include probe A32 to adapter A32 (!!!)
# Column A (pin index 1-32) is re-labeled "A32 (Probe Side)". A new
# "A32 (Adapter Side)" column is inserted as column B, and the original
# "Intan" values move to column C with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
# Column A keeps the original pin index, but is re-labeled as the probe-side
# pin. A new "A32 (Adapter Side)" column is inserted as column B, and the
# original "Intan" mapping moves to column C.
# (Order matters for shared-string table layout: adapter-side label first,
# then probe-side label, then the pre-existing "Intan" label last.)
$ws.Range("B1").Value = "A32 (Adapter Side)"
$ws.Range("A1").Value = "A32 (Probe Side)"
$ws.Range("C1").Value = "Intan"

# --- New column B: A32 (Adapter Side) --------------------------------------
$adapterSide = @(16,6,5,15,4,7,3,8,2,9,1,10,14,13,12,11,22,21,20,19,23,25,24,18,26,17,27,29,28,31,30,32)

# --- New column C: Intan pin (re-derived so probe A32 lines up with adapter A32) ---
$intan = @(30,26,21,17,27,22,20,25,28,23,19,24,29,18,31,16,0,15,2,13,8,9,7,1,6,14,10,11,5,12,4,3)

for ($i = 0; $i -lt 32; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $adapterSide[$i]
    $ws.Cells.Item($row, 3).Value = $intan[$i]
}

# --- Column widths ----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21
$ws.Columns.Item(2).ColumnWidth = 23.5703125
$ws.Columns.Item(3).ColumnWidth = 18

# --- Selection ---------------------------------------------------------------
$ws.Range("G6").Select()

$wb.Save()
